# Updated cryptos list on Mon Mar 27 02:36:00 UTC 2023 with GitHub Actions
# Applies the latest coinranking.com snapshot: refreshed Price/Volume(1h) figures
# for existing rows, plus the OKB/Dogecoin (rows 9-10) and Filecoin/HuobiToken
# (rows 33-34) ranking swap that came with this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new literal text value. NumberFormat is forced to
# "@" (Text) before the write so numeric-looking strings (e.g. "1.003",
# "27.903.79") are stored verbatim instead of being parsed into floats, then
# restored to the workbook's default "Normal" style so no stray formatting
# is introduced.
$updates = [ordered]@{
    "D2" = "27.903.79"
    "E2" = "  +0.95%  "
    "D3" = "1.767.33"
    "E3" = "  +0.69%  "
    "D4" = "1.003"
    "E4" = "  +0.09%  "
    "D5" = "329.35"
    "E5" = "  +1.59%  "
    "E6" = "  +0.09%  "
    "D7" = "0.4470"
    "E7" = "  -0.52%  "
    "D8" = "0.3521"
    "E8" = "  -0.81%  "
    "B9" = "OKB"
    "C9" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D9" = "41.89"
    "E9" = "  +0.66%  "
    "B10" = "Dogecoin"
    "C10" = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
    "D10" = "0.07371"
    "E10" = "  -1.16%  "
    "D11" = "1.093"
    "E11" = "  +0.46%  "
    "E12" = "  +0.09%  "
    "D13" = "20.67"
    "E13" = "  -0.32%  "
    "D14" = "6.007"
    "E14" = "  +0.20%  "
    "D15" = "7.187"
    "E15" = "  +0.23%  "
    "D16" = "1.770.81"
    "E16" = "  +1.00%  "
    "D17" = "92.42"
    "E17" = "  -0.98%  "
    "D18" = "0.00001059"
    "E18" = "  -0.22%  "
    "D19" = "0.06421"
    "E19" = "  -0.70%  "
    "E20" = "  +0.11%  "
    "D21" = "16.95"
    "E21" = "  -0.62%  "
    "D22" = "5.768"
    "E22" = "  +0.36%  "
    "D23" = "27.933.42"
    "E23" = "  +0.91%  "
    "D24" = "11.20"
    "E24" = "  -0.30%  "
    "D25" = "2.103"
    "E25" = "  -0.54%  "
    "D26" = "160.15"
    "E26" = "  -2.40%  "
    "D27" = "20.15"
    "E27" = "  -0.05%  "
    "D28" = "1.974.01"
    "E28" = "  +0.95%  "
    "D29" = "2.127"
    "E29" = "  +2.01%  "
    "D30" = "124.13"
    "E30" = "  -0.97%  "
    "D31" = "1.082"
    "E31" = "  +1.05%  "
    "D32" = "0.09146"
    "E32" = "  -0.32%  "
    "B33" = "HuobiToken"
    "C33" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "D33" = "3.676"
    "E33" = "  +0.56%  "
    "B34" = "Filecoin"
    "C34" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
    "D34" = "5.600"
    "E34" = "  +2.09%  "
    "D35" = "11.81"
    "E35" = "  +0.88%  "
    "D36" = "0.02279"
    "E36" = "  -0.23%  "
    "D37" = "0.06091"
    "E37" = "  +0.74%  "
    "D38" = "0.2084"
    "E38" = "  +0.45%  "
    "D39" = "4.949"
    "E39" = "  -0.09%  "
    "D40" = "0.6245"
    "E40" = "  -0.82%  "
    "D41" = "1.178"
    "E41" = "  -0.36%  "
    "D42" = "1.383"
    "E42" = "  -0.35%  "
    "D43" = "7.789"
    "E43" = "  +0.29%  "
    "D44" = "13.26"
    "E44" = "  +0.58%  "
    "D45" = "3.738"
    "E45" = "  +0.75%  "
    "D46" = "0.5821"
    "E46" = "  -1.03%  "
    "D47" = "122.16"
    "E47" = "  -0.82%  "
    "D48" = "1.933"
    "E48" = "  -0.59%  "
    "D49" = "1.133"
    "E49" = "  +0.47%  "
    "D50" = "0.06838"
    "E50" = "  -0.87%  "
    "E51" = "  +1.74%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}

